$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("employees")

# New employee record - row 11 (employeeid 1017, Monika, backend)
$ws.Range("A11").Value = 1017
$ws.Range("B11").Value = "Monika"
$ws.Range("C11").Value = "Mgonila26@gmail.com"
$ws.Range("D11").Value = "backend"
$ws.Range("E11").Value = "http://res.cloudinary.com/db5vn6bj5/image/upload/v1631979713/pdakv5si8r3yiukxa5ht.jpg"
$ws.Range("F11").Value = "experience"
$ws.Range("G11").Value = "delhi"
$ws.Range("H11").Value = "B.tech"
$ws.Range("I11").Value = "full time"
$ws.Range("J11").Value = 44485.44639960648
$ws.Range("J11").NumberFormat = "m/d/yy"
$ws.Range("K11").Value = 9140834289

# New employee record - row 12 (employeeid 1018, shreyansh, backend)
$ws.Range("A12").Value = 1018
$ws.Range("B12").Value = "shreyansh"
$ws.Range("C12").Value = "sj26@gmail.com"
$ws.Range("D12").Value = "backend"
$ws.Range("E12").Value = "http://res.cloudinary.com/db5vn6bj5/image/upload/v1631979713/pdakv5si8r3yiukxa5ht.jpg"
$ws.Range("F12").Value = "experience"
$ws.Range("G12").Value = "delhi"
$ws.Range("H12").Value = "B.tech"
$ws.Range("I12").Value = "full time"
$ws.Range("J12").Value = 44486.421475219904
$ws.Range("J12").NumberFormat = "m/d/yy"
$ws.Range("K12").Value = 9140834289
